$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: fill with same sample data as row 4 (Ali / Lopez Sarmiento / QC / Collège des médecins du Québec / 4426 / VERIFIED)
$ws.Range("A5").Value = "Ali"
$ws.Range("B5").Value = "Lopez Sarmiento"
$ws.Range("C5").Value = "QC"
$ws.Range("D5").Value = "Collège des médecins du Québec"
$ws.Range("E5").Value = 4426
$ws.Range("F5").Value = "VERIFIED"

# Row 6: same sample data copied down again
$ws.Range("A6").Value = "Ali"
$ws.Range("B6").Value = "Lopez Sarmiento"
$ws.Range("C6").Value = "QC"
$ws.Range("D6").Value = "Collège des médecins du Québec"
$ws.Range("E6").Value = 4426
$ws.Range("F6").Value = "VERIFIED"

# D6 previously carried an explicit (bold) style with no content; clear it back to Normal
# now that it holds a value, matching the formatting already used by D4/D5.
$ws.Range("D6").Style = "Normal"

# Update the active selection to D11
$ws.Range("D11").Select()
